$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.024436122217634
$ws.Cells.Item(2, 4).Value = 1.030658426276773
$ws.Cells.Item(2, 5).Value = 1.024925770673861
$ws.Cells.Item(2, 6).Value = 1.038684624685416
$ws.Cells.Item(2, 9).Value = 1.034742927218476
$ws.Cells.Item(2, 10).Value = 1.029610747642374
$ws.Cells.Item(2, 11).Value = 1.033468799841401
$ws.Cells.Item(2, 12).Value = 1.027752842327732
$ws.Cells.Item(2, 13).Value = 1.041471955892402
$ws.Cells.Item(2, 14).Value = 1.01391865233304
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.025223130670642
$ws.Cells.Item(3, 4).Value = 1.031275793302497
$ws.Cells.Item(3, 5).Value = 1.025589532947276
$ws.Cells.Item(3, 6).Value = 1.040867152547867
$ws.Cells.Item(3, 9).Value = 1.035001303644934
$ws.Cells.Item(3, 10).Value = 1.030037502238619
$ws.Cells.Item(3, 11).Value = 1.033895080878977
$ws.Cells.Item(3, 12).Value = 1.028224189043853
$ws.Cells.Item(3, 13).Value = 1.043460928340301
$ws.Cells.Item(3, 14).Value = 1.014060813111157
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.025732602309017
$ws.Cells.Item(4, 4).Value = 1.031675256021964
$ws.Cells.Item(4, 5).Value = 1.026019616043011
$ws.Cells.Item(4, 6).Value = 1.042273965448634
$ws.Cells.Item(4, 9).Value = 1.035166857675854
$ws.Cells.Item(4, 10).Value = 1.030313200750431
$ws.Cells.Item(4, 11).Value = 1.034170206978339
$ws.Cells.Item(4, 12).Value = 1.028529084917112
$ws.Cells.Item(4, 13).Value = 1.044742166789552
$ws.Cells.Item(4, 14).Value = 1.014152627259865
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.02594683753615
$ws.Cells.Item(5, 4).Value = 1.031843186370521
$ws.Cells.Item(5, 5).Value = 1.026200562195709
$ws.Cells.Item(5, 6).Value = 1.042864121164281
$ws.Cells.Item(5, 9).Value = 1.035236066978097
$ws.Cells.Item(5, 10).Value = 1.03042899901916
$ws.Cells.Item(5, 11).Value = 1.034285701186481
$ws.Cells.Item(5, 12).Value = 1.028657239666248
$ws.Cells.Item(5, 13).Value = 1.045279449562982
$ws.Cells.Item(5, 14).Value = 1.014191184385633
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.025982811689776
$ws.Cells.Item(6, 4).Value = 1.031871382387588
$ws.Cells.Item(6, 5).Value = 1.026230951993605
$ws.Cells.Item(6, 6).Value = 1.042963137422625
$ws.Cells.Item(6, 9).Value = 1.035247664714099
$ws.Cells.Item(6, 10).Value = 1.030448435878077
$ws.Cells.Item(6, 11).Value = 1.03430508328086
$ws.Cells.Item(6, 12).Value = 1.028678756027938
$ws.Cells.Item(6, 13).Value = 1.045369583383425
$ws.Cells.Item(6, 14).Value = 1.014197655859507
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.025735464721528
$ws.Cells.Item(7, 4).Value = 1.031677499930639
$ws.Cells.Item(7, 5).Value = 1.026022033308269
$ws.Cells.Item(7, 6).Value = 1.042281856084327
$ws.Cells.Item(7, 9).Value = 1.035167783983312
$ws.Cells.Item(7, 10).Value = 1.03031474846669
$ws.Cells.Item(7, 11).Value = 1.034171750880603
$ws.Cells.Item(7, 12).Value = 1.028530797419756
$ws.Cells.Item(7, 13).Value = 1.044749351257081
$ws.Cells.Item(7, 14).Value = 1.014153142625358
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.024702048678381
$ws.Cells.Item(8, 4).Value = 1.030867071365438
$ws.Cells.Item(8, 5).Value = 1.025149970789172
$ws.Cells.Item(8, 6).Value = 1.039423363106285
$ws.Cells.Item(8, 9).Value = 1.034830585603
$ws.Cells.Item(8, 10).Value = 1.02975506252535
$ws.Cells.Item(8, 11).Value = 1.03361300991539
$ws.Cells.Item(8, 12).Value = 1.027912156084555
$ws.Cells.Item(8, 13).Value = 1.042145348613205
$ws.Cells.Item(8, 14).Value = 1.013966732114564
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.02288277375731
$ws.Cells.Item(9, 4).Value = 1.029438883694142
$ws.Cells.Item(9, 5).Value = 1.023617792242078
$ws.Cells.Item(9, 6).Value = 1.034343344500414
$ws.Cells.Item(9, 9).Value = 1.034223834714504
$ws.Cells.Item(9, 10).Value = 1.02876545072028
$ws.Cells.Item(9, 11).Value = 1.032623014284017
$ws.Cells.Item(9, 12).Value = 1.026821298596753
$ws.Cells.Item(9, 13).Value = 1.037511335988925
$ws.Cells.Item(9, 14).Value = 1.013636926146568
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.021671109521123
$ws.Cells.Item(10, 4).Value = 1.028486689260114
$ws.Cells.Item(10, 5).Value = 1.022599408113143
$ws.Cells.Item(10, 6).Value = 1.030925770501917
$ws.Cells.Item(10, 9).Value = 1.033810801195839
$ws.Cells.Item(10, 10).Value = 1.028103435332792
$ws.Cells.Item(10, 11).Value = 1.031959348551822
$ws.Cells.Item(10, 12).Value = 1.026093578748841
$ws.Cells.Item(10, 13).Value = 1.034389597206827
$ws.Cells.Item(10, 14).Value = 1.013416163744083
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.021146729445341
$ws.Cells.Item(11, 4).Value = 1.028074361602587
$ws.Cells.Item(11, 5).Value = 1.022159171388318
$ws.Cells.Item(11, 6).Value = 1.02943814729716
$ws.Cells.Item(11, 9).Value = 1.033629909320283
$ws.Cells.Item(11, 10).Value = 1.027816233632652
$ws.Cells.Item(11, 11).Value = 1.031671097658915
$ws.Cells.Item(11, 12).Value = 1.025778355184925
$ws.Cells.Item(11, 13).Value = 1.03302974841749
$ws.Cells.Item(11, 14).Value = 1.013320359393315
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.020951993282463
$ws.Cells.Item(12, 4).Value = 1.027921201399793
$ws.Cells.Item(12, 5).Value = 1.021995758084274
$ws.Cells.Item(12, 6).Value = 1.028884368371682
$ws.Cells.Item(12, 9).Value = 1.033562408883147
$ws.Cells.Item(12, 10).Value = 1.027709472007467
$ws.Cells.Item(12, 11).Value = 1.0315638957116
$ws.Cells.Item(12, 12).Value = 1.025661249701516
$ws.Cells.Item(12, 13).Value = 1.032523384548017
$ws.Cells.Item(12, 14).Value = 1.013284741366874
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.020993762951065
$ws.Cells.Item(13, 4).Value = 1.027954054938022
$ws.Cells.Item(13, 5).Value = 1.022030805797532
$ws.Cells.Item(13, 6).Value = 1.029003211191982
$ws.Cells.Item(13, 9).Value = 1.033576901974782
$ws.Cells.Item(13, 10).Value = 1.027732376459911
$ws.Cells.Item(13, 11).Value = 1.031586896910129
$ws.Cells.Item(13, 12).Value = 1.025686370013301
$ws.Cells.Item(13, 13).Value = 1.032632058718683
$ws.Cells.Item(13, 14).Value = 1.013292383004681
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.021130631637237
$ws.Cells.Item(14, 4).Value = 1.028061701401948
$ws.Cells.Item(14, 5).Value = 1.022145661342484
$ws.Cells.Item(14, 6).Value = 1.029392396657937
$ws.Cells.Item(14, 9).Value = 1.033624336030077
$ws.Cells.Item(14, 10).Value = 1.027807410360574
$ws.Cells.Item(14, 11).Value = 1.031662239023984
$ws.Cells.Item(14, 12).Value = 1.0257686755644
$ws.Cells.Item(14, 13).Value = 1.032987918045505
$ws.Cells.Item(14, 14).Value = 1.013317415849808
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.021214966512617
$ws.Cells.Item(15, 4).Value = 1.028128025491241
$ws.Cells.Item(15, 5).Value = 1.022216442253541
$ws.Cells.Item(15, 6).Value = 1.029632025231492
$ws.Cells.Item(15, 9).Value = 1.033653520698141
$ws.Cells.Item(15, 10).Value = 1.02785363033135
$ws.Cells.Item(15, 11).Value = 1.031708642179165
$ws.Cells.Item(15, 12).Value = 1.025819384423184
$ws.Cells.Item(15, 13).Value = 1.033207007262775
$ws.Cells.Item(15, 14).Value = 1.013332835166366
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.021705916777247
$ws.Cells.Item(16, 4).Value = 1.028514053628916
$ws.Cells.Item(16, 5).Value = 1.02262864060463
$ws.Cells.Item(16, 6).Value = 1.031024331699396
$ws.Cells.Item(16, 9).Value = 1.033822763141487
$ws.Cells.Item(16, 10).Value = 1.028122484463622
$ws.Cells.Item(16, 11).Value = 1.031978460236163
$ws.Cells.Item(16, 12).Value = 1.026114496655645
$ws.Cells.Item(16, 13).Value = 1.03447967183252
$ws.Cells.Item(16, 14).Value = 1.013422517475691
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.022013951557583
$ws.Cells.Item(17, 4).Value = 1.028756193102047
$ws.Cells.Item(17, 5).Value = 1.022887397626371
$ws.Cells.Item(17, 6).Value = 1.031895575605389
$ws.Cells.Item(17, 9).Value = 1.033928375518562
$ws.Cells.Item(17, 10).Value = 1.028290983534459
$ws.Cells.Item(17, 11).Value = 1.032147474137868
$ws.Cells.Item(17, 12).Value = 1.026299581711
$ws.Cells.Item(17, 13).Value = 1.035275782603245
$ws.Cells.Item(17, 14).Value = 1.013478715815398
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.02219364990415
$ws.Cells.Item(18, 4).Value = 1.028897426930725
$ws.Cells.Item(18, 5).Value = 1.023038396654319
$ws.Cells.Item(18, 6).Value = 1.032403008367044
$ws.Cells.Item(18, 9).Value = 1.033989780181886
$ws.Cells.Item(18, 10).Value = 1.028389213661352
$ws.Cells.Item(18, 11).Value = 1.03224597229772
$ws.Cells.Item(18, 12).Value = 1.026407527494406
$ws.Cells.Item(18, 13).Value = 1.035739360017735
$ws.Cells.Item(18, 14).Value = 1.013511474831285
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.022254926950333
$ws.Cells.Item(19, 4).Value = 1.028945583676405
$ws.Cells.Item(19, 5).Value = 1.023089895351847
$ws.Cells.Item(19, 6).Value = 1.032575903938023
$ws.Cells.Item(19, 9).Value = 1.034010684182126
$ws.Cells.Item(19, 10).Value = 1.028422698694295
$ws.Cells.Item(19, 11).Value = 1.032279543261378
$ws.Cells.Item(19, 12).Value = 1.026444332316235
$ws.Cells.Item(19, 13).Value = 1.035897296632688
$ws.Cells.Item(19, 14).Value = 1.013522641329268
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.02198089956089
$ws.Cells.Item(20, 4).Value = 1.028730214039326
$ws.Cells.Item(20, 5).Value = 1.022859628157862
$ws.Cells.Item(20, 6).Value = 1.031802177138969
$ws.Cells.Item(20, 9).Value = 1.033917064725456
$ws.Cells.Item(20, 10).Value = 1.028272910614445
$ws.Cells.Item(20, 11).Value = 1.032129349320518
$ws.Cells.Item(20, 12).Value = 1.026279724994333
$ws.Cells.Item(20, 13).Value = 1.035190448446215
$ws.Cells.Item(20, 14).Value = 1.013472688388816
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.021090326061981
$ws.Cells.Item(21, 4).Value = 1.028030002295424
$ws.Cells.Item(21, 5).Value = 1.022111836200842
$ws.Cells.Item(21, 6).Value = 1.029277824940364
$ws.Cells.Item(21, 9).Value = 1.033610376431575
$ws.Cells.Item(21, 10).Value = 1.027785317016599
$ws.Cells.Item(21, 11).Value = 1.031640056318608
$ws.Cells.Item(21, 12).Value = 1.025744439111754
$ws.Cells.Item(21, 13).Value = 1.032883161294978
$ws.Cells.Item(21, 14).Value = 1.013310045185402
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.02053062986555
$ws.Cells.Item(22, 4).Value = 1.027589731874261
$ws.Cells.Item(22, 5).Value = 1.02164230762229
$ws.Cells.Item(22, 6).Value = 1.027683647367668
$ws.Cells.Item(22, 9).Value = 1.033415759809929
$ws.Cells.Item(22, 10).Value = 1.027478272244911
$ws.Cells.Item(22, 11).Value = 1.031331650105365
$ws.Cells.Item(22, 12).Value = 1.025407782856335
$ws.Cells.Item(22, 13).Value = 1.031425196042087
$ws.Cells.Item(22, 14).Value = 1.013207599626071
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.020827312377478
$ws.Cells.Item(23, 4).Value = 1.027823129470303
$ws.Cells.Item(23, 5).Value = 1.021891152982038
$ws.Cells.Item(23, 6).Value = 1.028529429371392
$ws.Cells.Item(23, 9).Value = 1.033519100003903
$ws.Cells.Item(23, 10).Value = 1.02764108764916
$ws.Cells.Item(23, 11).Value = 1.031495215119332
$ws.Cells.Item(23, 12).Value = 1.025586260259089
$ws.Cells.Item(23, 13).Value = 1.032198793547896
$ws.Cells.Item(23, 14).Value = 1.013261925543643
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.02199583425291
$ws.Cells.Item(24, 4).Value = 1.028741952866127
$ws.Cells.Item(24, 5).Value = 1.022872175768113
$ws.Cells.Item(24, 6).Value = 1.031844382198761
$ws.Cells.Item(24, 9).Value = 1.033922176195565
$ws.Cells.Item(24, 10).Value = 1.028281077151618
$ws.Cells.Item(24, 11).Value = 1.032137539407527
$ws.Cells.Item(24, 12).Value = 1.026288697424751
$ws.Cells.Item(24, 13).Value = 1.035229009689034
$ws.Cells.Item(24, 14).Value = 1.013475411986946
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.023352892192624
$ws.Cells.Item(25, 4).Value = 1.029808116818084
$ws.Cells.Item(25, 5).Value = 1.024013359006737
$ws.Cells.Item(25, 6).Value = 1.03566194853988
$ws.Cells.Item(25, 9).Value = 1.034382191946740
$ws.Cells.Item(25, 10).Value = 1.029021689405222
$ws.Cells.Item(25, 11).Value = 1.03287959672488
$ws.Cells.Item(25, 12).Value = 1.027103396987384
$ws.Cells.Item(25, 13).Value = 1.038714911743198
$ws.Cells.Item(25, 14).Value = 1.013722346065132
